$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header row: C1/D1 should use the same (bold Cambria) style as B1 ---
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)  # xlPasteFormats

# --- Append the new CA_Market table below the existing one (rows 16-20) ---
# Row 16: header row, same layout/style as row 1
$ws.Range("B1:D1").Copy($ws.Range("B16:D16"))

# Row 17: precision (copy format/labels from row 2, then overwrite values)
$ws.Range("A2:D2").Copy($ws.Range("A17:D17"))
$ws.Range("B17").Value2 = 0.5813
$ws.Range("C17").Value2 = 0.42
$ws.Range("D17").Value2 = 0.3277

# Row 18: recall
$ws.Range("A3:D3").Copy($ws.Range("A18:D18"))
$ws.Range("B18").Value2 = 0.5413
$ws.Range("C18").Value2 = 0.3934
$ws.Range("D18").Value2 = 0.328

# Row 19: f1
$ws.Range("A5:D5").Copy($ws.Range("A19:D19"))
$ws.Range("B19").Value2 = 0.5434
$ws.Range("C19").Value2 = 0.3816
$ws.Range("D19").Value2 = 0.311

# Row 20: MA
$ws.Range("A6:D6").Copy($ws.Range("A20:D20"))
$ws.Range("B20").Value2 = 0.71688574552536
$ws.Range("C20").Value2 = 0.634228229522705
$ws.Range("D20").Value2 = 0.602592468261719

# --- Move the active selection to A20, matching the author's last edit point ---
$ws.Range("A20").Select()
